# "Add files via upload" update:
#   - Bump the confidential-notice "as of" date from 2021-03-25 to 2021-03-26
#   - Refresh the Weight / Percent Change values in D2:E8
#
# The sheet ships protected, so it has to be unprotected for the duration of
# the edit and re-protected afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential notice text (cell A11) with the new date.
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-8.
$ws.Range("D2").Value = 0.4995226079454501
$ws.Range("E2").Value = 0.01607372481783109

$ws.Range("D3").Value = 0.2423184891337296
$ws.Range("E3").Value = 0.01588799748308944

$ws.Range("D4").Value = 0.09782106871343398
$ws.Range("E4").Value = 0.02468007312614251

$ws.Range("D5").Value = 0.102118712369513
$ws.Range("E5").Value = 0.02084781097984711

$ws.Range("D6").Value = 0.03024631970162236
$ws.Range("E6").Value = 0.02413273001508287

$ws.Range("D7").Value = 0.02797280213625088
$ws.Range("E7").Value = 0.02280338064104592

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = 0.01779012794999191

# Re-apply sheet protection.
$ws.Protect()
